$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -------------------------------------------------
# Insert a new blank column B ("timestamp"); annotator .. quotes shift right.
$ws.Range("B1").EntireColumn.Insert()

# Remove the old "periodical" column (now shifted to column L).
$ws.Range("L1").EntireColumn.Delete()

# Insert the 3 new data rows at their (alphabetically) sorted positions.
# Doing this top-down with the FINAL row numbers works because each insert
# only pushes rows at/after it down by one, and we always target the next
# still-to-be-placed new row in top-to-bottom order.
$ws.Range("A3").EntireRow.Insert()
$ws.Range("A4").EntireRow.Insert()
$ws.Range("A9").EntireRow.Insert()

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1,1).Value = "file"
$ws.Cells.Item(1,2).Value = "timestamp"
$ws.Cells.Item(1,3).Value = "annotator"
$ws.Cells.Item(1,4).Value = "orig_filename"
$ws.Cells.Item(1,5).Value = "year"
$ws.Cells.Item(1,6).Value = "decade"
$ws.Cells.Item(1,7).Value = "source"
$ws.Cells.Item(1,8).Value = "title"
$ws.Cells.Item(1,9).Value = "author"
$ws.Cells.Item(1,10).Value = "fictional"
$ws.Cells.Item(1,11).Value = "text_type"
$ws.Cells.Item(1,12).Value = "narrative"
$ws.Cells.Item(1,13).Value = "cabtokens"
$ws.Cells.Item(1,14).Value = "dialect"
$ws.Cells.Item(1,15).Value = "perspective"
$ws.Cells.Item(1,16).Value = "quotes"

# --- Force year (E) and decade (F) columns to be stored as text, matching
# the workbook author's original authoring (these look like numbers but are
# typed as text in the source file).
$ws.Range("E2:F16").NumberFormat = "@"

# --- Data rows --------------------------------------------------------
# row 2: rwz_full_grenz_10479.xmi
$ws.Cells.Item(2,1).Value = "rwz_full_grenz_10479.xmi"
$ws.Cells.Item(2,3).Value = "mercury"
$ws.Cells.Item(2,4).Value = "grenzboten_341831_157679.txt_10.xml"
$ws.Cells.Item(2,5).Value = "1880"
$ws.Cells.Item(2,6).Value = "1880"
$ws.Cells.Item(2,7).Value = "grenz"
$ws.Cells.Item(2,8).Value = "Die Tragik in Werken der hellenischen Plastik"
$ws.Cells.Item(2,9).Value = "Veit, Valentin"
$ws.Cells.Item(2,10).Value = "no"
$ws.Cells.Item(2,11).Value = "Undefined"
$ws.Cells.Item(2,12).Value = "no"
$ws.Cells.Item(2,13).Value = 6609
$ws.Cells.Item(2,14).Value = "NA"
$ws.Cells.Item(2,15).Value = "NA"
$ws.Cells.Item(2,16).Value = "NA"

# row 3: rwz_full_grenz_13.xmi
$ws.Cells.Item(3,1).Value = "rwz_full_grenz_13.xmi"
$ws.Cells.Item(3,3).Value = "uranus"
$ws.Cells.Item(3,4).Value = "grenzboten_282160_266616.txt_20.xml"
$ws.Cells.Item(3,5).Value = "1842"
$ws.Cells.Item(3,6).Value = "1840"
$ws.Cells.Item(3,7).Value = "grenz"
$ws.Cells.Item(3,8).Value = "Der Musikunterricht in Elementarschulen in Deutschland und Frankreich."
$ws.Cells.Item(3,9).Value = "Undefined"
$ws.Cells.Item(3,10).Value = "no"
$ws.Cells.Item(3,11).Value = "undefined"
$ws.Cells.Item(3,12).Value = "no"
$ws.Cells.Item(3,13).Value = 9240
$ws.Cells.Item(3,14).Value = "NA"
$ws.Cells.Item(3,15).Value = "NA"
$ws.Cells.Item(3,16).Value = "NA"

# row 4: rwz_full_grenz_18399.xmi
$ws.Cells.Item(4,1).Value = "rwz_full_grenz_18399.xmi"
$ws.Cells.Item(4,3).Value = "uranus"
$ws.Cells.Item(4,4).Value = "grenzboten_341875_238787.txt_44.xml"
$ws.Cells.Item(4,5).Value = "1902"
$ws.Cells.Item(4,6).Value = "1900"
$ws.Cells.Item(4,7).Value = "grenz"
$ws.Cells.Item(4,8).Value = "Musikalische Zeitfragen"
$ws.Cells.Item(4,9).Value = "Kretzschmar, Hermann"
$ws.Cells.Item(4,10).Value = "no"
$ws.Cells.Item(4,11).Value = "Undefined"
$ws.Cells.Item(4,12).Value = "no"
$ws.Cells.Item(4,13).Value = 7195
$ws.Cells.Item(4,14).Value = "NA"
$ws.Cells.Item(4,15).Value = "NA"
$ws.Cells.Item(4,16).Value = "NA"

# row 5: rwz_full_grenz_20389.xmi
$ws.Cells.Item(5,1).Value = "rwz_full_grenz_20389.xmi"
$ws.Cells.Item(5,3).Value = "neptune"
$ws.Cells.Item(5,4).Value = "grenzboten_341885_303415.txt_6.xml"
$ws.Cells.Item(5,5).Value = "1907"
$ws.Cells.Item(5,6).Value = "1900"
$ws.Cells.Item(5,7).Value = "grenz"
$ws.Cells.Item(5,8).Value = "Goethe und die Boisserée"
$ws.Cells.Item(5,9).Value = "Undefined"
$ws.Cells.Item(5,10).Value = "no"
$ws.Cells.Item(5,11).Value = "Undefined"
$ws.Cells.Item(5,12).Value = "no"
$ws.Cells.Item(5,13).Value = 4483
$ws.Cells.Item(5,14).Value = "NA"
$ws.Cells.Item(5,15).Value = "NA"
$ws.Cells.Item(5,16).Value = "NA"

# row 6: rwz_full_grenz_23110.xmi
$ws.Cells.Item(6,1).Value = "rwz_full_grenz_23110.xmi"
$ws.Cells.Item(6,3).Value = "mercury"
$ws.Cells.Item(6,4).Value = "grenzboten_341901_323097.txt_30.xml"
$ws.Cells.Item(6,5).Value = "1915"
$ws.Cells.Item(6,6).Value = "1910"
$ws.Cells.Item(6,7).Value = "grenz"
$ws.Cells.Item(6,8).Value = "Die Stellung Belgiens zum alten Reiche"
$ws.Cells.Item(6,9).Value = "Born, Conrad Professor Dr."
$ws.Cells.Item(6,10).Value = "no"
$ws.Cells.Item(6,11).Value = "Undefined"
$ws.Cells.Item(6,12).Value = "no"
$ws.Cells.Item(6,13).Value = 5091
$ws.Cells.Item(6,14).Value = "NA"
$ws.Cells.Item(6,15).Value = "NA"
$ws.Cells.Item(6,16).Value = "NA"

# row 7: rwz_full_grenz_3721.xmi
$ws.Cells.Item(7,1).Value = "rwz_full_grenz_3721.xmi"
$ws.Cells.Item(7,3).Value = "venus"
$ws.Cells.Item(7,4).Value = "grenzboten_341584_102594.txt_38.xml"
$ws.Cells.Item(7,5).Value = "1856"
$ws.Cells.Item(7,6).Value = "1850"
$ws.Cells.Item(7,7).Value = "grenz"
$ws.Cells.Item(7,8).Value = "Regierung und Volk in Neapel"
$ws.Cells.Item(7,9).Value = "Undefined"
$ws.Cells.Item(7,10).Value = "no"
$ws.Cells.Item(7,11).Value = "Undefined"
$ws.Cells.Item(7,12).Value = "no"
$ws.Cells.Item(7,13).Value = 8039
$ws.Cells.Item(7,14).Value = "NA"
$ws.Cells.Item(7,15).Value = "NA"
$ws.Cells.Item(7,16).Value = "NA"

# row 8: rwz_full_grenz_6595.xmi
$ws.Cells.Item(8,1).Value = "rwz_full_grenz_6595.xmi"
$ws.Cells.Item(8,3).Value = "neptune"
$ws.Cells.Item(8,4).Value = "grenzboten_341805_191229.txt_24.xml"
$ws.Cells.Item(8,5).Value = "1867"
$ws.Cells.Item(8,6).Value = "1860"
$ws.Cells.Item(8,7).Value = "grenz"
$ws.Cells.Item(8,8).Value = "Die Petrussagen"
$ws.Cells.Item(8,9).Value = "Undefined"
$ws.Cells.Item(8,10).Value = "no"
$ws.Cells.Item(8,11).Value = "Undefined"
$ws.Cells.Item(8,12).Value = "no"
$ws.Cells.Item(8,13).Value = 5833
$ws.Cells.Item(8,14).Value = "NA"
$ws.Cells.Item(8,15).Value = "NA"
$ws.Cells.Item(8,16).Value = "NA"

# row 9: rwz_full_grenz_7300.xmi
$ws.Cells.Item(9,1).Value = "rwz_full_grenz_7300.xmi"
$ws.Cells.Item(9,3).Value = "uranus"
$ws.Cells.Item(9,4).Value = "grenzboten_341809_121220.txt_26.xml"
$ws.Cells.Item(9,5).Value = "1869"
$ws.Cells.Item(9,6).Value = "1860"
$ws.Cells.Item(9,7).Value = "grenz"
$ws.Cells.Item(9,8).Value = "Polnischer Monatsbericht."
$ws.Cells.Item(9,9).Value = "Undefined"
$ws.Cells.Item(9,10).Value = "no"
$ws.Cells.Item(9,11).Value = "Undefined"
$ws.Cells.Item(9,12).Value = "no"
$ws.Cells.Item(9,13).Value = 7071
$ws.Cells.Item(9,14).Value = "NA"
$ws.Cells.Item(9,15).Value = "NA"
$ws.Cells.Item(9,16).Value = "NA"

# row 10: rwz_full_mkhz_10107.xmi
$ws.Cells.Item(10,1).Value = "rwz_full_mkhz_10107.xmi"
$ws.Cells.Item(10,3).Value = "neptune"
$ws.Cells.Item(10,4).Value = "MT_1887_01_10.txt_1.xml"
$ws.Cells.Item(10,5).Value = "1887"
$ws.Cells.Item(10,6).Value = "1880"
$ws.Cells.Item(10,7).Value = "mkhz.maehrisches"
$ws.Cells.Item(10,8).Value = "Die Versammlung der deutschen Vertrauensmänner in Brünn"
$ws.Cells.Item(10,9).Value = "Undefined"
$ws.Cells.Item(10,10).Value = "no"
$ws.Cells.Item(10,11).Value = "Undefined"
$ws.Cells.Item(10,12).Value = "no"
$ws.Cells.Item(10,13).Value = 7131
$ws.Cells.Item(10,14).Value = "NA"
$ws.Cells.Item(10,15).Value = "NA"
$ws.Cells.Item(10,16).Value = "NA"

# row 11: rwz_full_mkhz_10440.xmi
$ws.Cells.Item(11,1).Value = "rwz_full_mkhz_10440.xmi"
$ws.Cells.Item(11,3).Value = "neptune"
$ws.Cells.Item(11,4).Value = "MT_1895_05_24.txt_1.xml"
$ws.Cells.Item(11,5).Value = "1895"
$ws.Cells.Item(11,6).Value = "1890"
$ws.Cells.Item(11,7).Value = "mkhz.maehrisches"
$ws.Cells.Item(11,8).Value = "Die 15. Hauptversammlung des Deutschen Schulvereins"
$ws.Cells.Item(11,9).Value = "Undefined"
$ws.Cells.Item(11,10).Value = "no"
$ws.Cells.Item(11,11).Value = "Undefined"
$ws.Cells.Item(11,12).Value = "no"
$ws.Cells.Item(11,13).Value = 6257
$ws.Cells.Item(11,14).Value = "NA"
$ws.Cells.Item(11,15).Value = "NA"
$ws.Cells.Item(11,16).Value = "NA"

# row 12: rwz_full_mkhz_1098.xmi
$ws.Cells.Item(12,1).Value = "rwz_full_mkhz_1098.xmi"
$ws.Cells.Item(12,3).Value = "mercury"
$ws.Cells.Item(12,4).Value = "nn_auswandererzeitung046_1852_2.xml"
$ws.Cells.Item(12,5).Value = "1852"
$ws.Cells.Item(12,6).Value = "1850"
$ws.Cells.Item(12,7).Value = "mkhz.auswandererzeitung"
$ws.Cells.Item(12,8).Value = "Die deutsche Kolonie Santa Cruz in der Provinz Rio Grande do Sul in Brasilien"
$ws.Cells.Item(12,9).Value = "Undefined"
$ws.Cells.Item(12,10).Value = "no"
$ws.Cells.Item(12,11).Value = "Undefined"
$ws.Cells.Item(12,12).Value = "no"
$ws.Cells.Item(12,13).Value = 3175
$ws.Cells.Item(12,14).Value = "NA"
$ws.Cells.Item(12,15).Value = "NA"
$ws.Cells.Item(12,16).Value = "NA"

# row 13: rwz_full_mkhz_2733.xmi
$ws.Cells.Item(13,1).Value = "rwz_full_mkhz_2733.xmi"
$ws.Cells.Item(13,3).Value = "mercury"
$ws.Cells.Item(13,4).Value = "nn_europa0106_1905_10.xml"
$ws.Cells.Item(13,5).Value = "1905"
$ws.Cells.Item(13,6).Value = "1900"
$ws.Cells.Item(13,7).Value = "mkhz.europa"
$ws.Cells.Item(13,8).Value = "Zur Beurteilung der jüngsten Schulkonflikte"
$ws.Cells.Item(13,9).Value = "Borchardt, Bruno"
$ws.Cells.Item(13,10).Value = "no"
$ws.Cells.Item(13,11).Value = "Undefined"
$ws.Cells.Item(13,12).Value = "no"
$ws.Cells.Item(13,13).Value = 3356
$ws.Cells.Item(13,14).Value = "NA"
$ws.Cells.Item(13,15).Value = "NA"
$ws.Cells.Item(13,16).Value = "NA"

# row 14: rwz_full_mkhz_336.xmi
$ws.Cells.Item(14,1).Value = "rwz_full_mkhz_336.xmi"
$ws.Cells.Item(14,3).Value = "neptune"
$ws.Cells.Item(14,4).Value = "nn_auswanderer16_1848.txt_3.xml"
$ws.Cells.Item(14,5).Value = "1848"
$ws.Cells.Item(14,6).Value = "1840"
$ws.Cells.Item(14,7).Value = "mkhz.auswanderer"
$ws.Cells.Item(14,8).Value = "Undefined"
$ws.Cells.Item(14,9).Value = "Undefined"
$ws.Cells.Item(14,10).Value = "no"
$ws.Cells.Item(14,11).Value = "Undefined"
$ws.Cells.Item(14,12).Value = "no"
$ws.Cells.Item(14,13).Value = 3368
$ws.Cells.Item(14,14).Value = "NA"
$ws.Cells.Item(14,15).Value = "NA"
$ws.Cells.Item(14,16).Value = "NA"

# row 15: rwz_full_mkhz_5884.xmi
$ws.Cells.Item(15,1).Value = "rwz_full_mkhz_5884.xmi"
$ws.Cells.Item(15,3).Value = "neptune"
$ws.Cells.Item(15,4).Value = "nn_social04_1873_11.xml"
$ws.Cells.Item(15,5).Value = "1873"
$ws.Cells.Item(15,6).Value = "1870"
$ws.Cells.Item(15,7).Value = "mkhz.social"
$ws.Cells.Item(15,8).Value = "Undefined"
$ws.Cells.Item(15,9).Value = "Undefined"
$ws.Cells.Item(15,10).Value = "no"
$ws.Cells.Item(15,11).Value = "Undefined"
$ws.Cells.Item(15,12).Value = "no"
$ws.Cells.Item(15,13).Value = 5118
$ws.Cells.Item(15,14).Value = "NA"
$ws.Cells.Item(15,15).Value = "NA"
$ws.Cells.Item(15,16).Value = "NA"

# row 16: rwz_full_mkhz_599.xmi
$ws.Cells.Item(16,1).Value = "rwz_full_mkhz_599.xmi"
$ws.Cells.Item(16,3).Value = "mercury"
$ws.Cells.Item(16,4).Value = "nn_auswanderer63_1847.txt_1.xml"
$ws.Cells.Item(16,5).Value = "1847"
$ws.Cells.Item(16,6).Value = "1840"
$ws.Cells.Item(16,7).Value = "mkhz.auswanderer"
$ws.Cells.Item(16,8).Value = "Undefined"
$ws.Cells.Item(16,9).Value = "Undefined"
$ws.Cells.Item(16,10).Value = "no"
$ws.Cells.Item(16,11).Value = "Undefined"
$ws.Cells.Item(16,12).Value = "no"
$ws.Cells.Item(16,13).Value = 4514
$ws.Cells.Item(16,14).Value = "NA"
$ws.Cells.Item(16,15).Value = "NA"
$ws.Cells.Item(16,16).Value = "NA"

# Remove the helper number-format now that the text values are committed,
# so cells end up with the default (no explicit) style, same as the rest
# of the sheet.
$ws.Range("E2:F16").ClearFormats()
